# The "South Korea" country block (rows 20-30) is being folded into the
# "United States" block above it: the standalone Country/City/Industry
# values for that block are cleared, its two "State1"/"State2" placeholder
# states are renamed to real US states (New Jersey, Texas), and column A's
# merge is extended so the whole A2:A30 range is one merged "United States"
# cell instead of two separate country cells (A2:A19 + A20:A30).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the now-unused Country cell (A20) and the City/Industry columns
# for rows 20-30 (they previously held placeholder City1/City2 + industry
# rows under the "South Korea" entry).
$ws.Range("A20").ClearContents()
$ws.Range("C20:D30").ClearContents()

# Rename the placeholder state names to real US states.
$ws.Range("B20").Value = "New Jersey"
$ws.Range("B27").Value = "Texas"

# Re-merge column A so United States spans the whole A2:A30 block instead
# of being split into two separate country merges.
$ws.Range("A20:A30").UnMerge()
$ws.Range("A2:A19").UnMerge()
$ws.Range("A2:A30").Merge()

# Update the active cell / selection on the sheet.
[void]$ws.Range("H12").Select()
